# =====================================================================
# Enhance job descriptions across the resume (Dheeraj Chand - Polling,
# Research & Redistricting / long / modern_tech variant).
#
# Strategy: for plain 1-for-1 text swaps, locate the old bullet with
# Find.Execute (no replacement string, so Word just selects/collapses
# the range onto the match) and then assign Range.Text directly. Doing
# the substitution this way - rather than passing the replacement
# through Find.Execute's ReplaceWith argument - avoids Word's
# AutoCorrect "smart quotes" silently turning a straight apostrophe
# into a curly one in the new text.
#
# For bullet lists whose paragraph counts grow, replace the bullets
# that have 1-for-1 counterparts the same way, then append the extra
# new bullets with Paragraph.Range.InsertParagraphAfter() + the
# newly created Paragraph's Range.Text.
# =====================================================================

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "WARNING: text not found -> $old"
        return
    }
    $r.Text = $new
}

function Add-BulletAfter($precedingText, $newBulletText) {
    $r = $d.Content
    $ok = $r.Find.Execute($precedingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "WARNING: anchor not found -> $precedingText"
        return $null
    }
    $p = $r.Paragraphs(1)
    $p.Range.InsertParagraphAfter()
    $newP = $p.Next()
    $newP.Range.Text = $newBulletText
    return $newP
}

# ---------------------------------------------------------------------
# PARTNER - Siege Analytics bullets
# ---------------------------------------------------------------------
Replace-Text "• Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations" "• Lead comprehensive polling and research studies for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions"

Replace-Text "• Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics" "• Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing millions of records with millions of columns for electoral analytics and demographic analysis"

Replace-Text "• Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets" "• Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets"

Replace-Text "• Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering" "• Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering"

Replace-Text "• Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications" "• Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications"

Replace-Text "• Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices" "• Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices"

# ---------------------------------------------------------------------
# DATA PRODUCTS MANAGER - Helm/Murmuration bullets
# ---------------------------------------------------------------------
Replace-Text "• Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES" "• Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES"

Replace-Text "• Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions" "• Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions"

Replace-Text "• Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI" "• Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI"

Replace-Text "• Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company's distinguishing products" "• Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company's distinguishing products"

Replace-Text "• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices" "• Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices"

# ---------------------------------------------------------------------
# SENIOR ANALYST - Myers Research bullet
# ---------------------------------------------------------------------
Replace-Text "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research" "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions"

# ---------------------------------------------------------------------
# RESEARCH DIRECTOR - Progressive Change Campaign Committee
# 4 bullets -> 5 bullets
# ---------------------------------------------------------------------
Replace-Text "• Managed critical research operations for political campaigns" "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls"

Replace-Text "• Conducted comprehensive polling and demographic analysis" "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren"

Replace-Text "• Developed strategic recommendations based on data analysis" "• Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"

Replace-Text "• Led research team in support of progressive political initiatives" "• Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly `$1 million annually in polling costs"

Add-BulletAfter "Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly `$1 million annually in polling costs" "• Managed comprehensive research operations for progressive political initiatives and candidates" | Out-Null

# ---------------------------------------------------------------------
# PROGRAMMER - Lake Research Partners
# heading + 4 bullets -> heading + 6 bullets
# (The section heading text "Political Research and Data Analysis" also
#  appears earlier for RESEARCH DIRECTOR, which must stay unchanged, so
#  scope this particular search to start right after the PROGRAMMER
#  job title.)
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("PROGRAMMER - Lake Research Partners, Washington, DC | April 2008", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$programmerScope = $d.Range($anchor.End, $d.Content.End)
$programmerScope.Find.Execute("Political Research and Data Analysis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$programmerScope.Text = "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns"

Replace-Text "• Developed data analysis tools for political polling and research" "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"

Replace-Text "• Built statistical models for voter behavior analysis" "• Developed system that later became the Polling Consortium Database at The Analyst Institute"

Replace-Text "• Created data visualization tools for research presentations" "• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions"

Replace-Text "• Supported senior researchers with technical analysis and reporting" "• Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle"

$p1 = Add-BulletAfter "Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle" "• Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps"
if ($p1) {
    $p1.Range.InsertParagraphAfter()
    $p2 = $p1.Next()
    $p2.Range.Text = "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding"
}

# ---------------------------------------------------------------------
# FIELD DIRECTOR - The Feldman Group
# heading + 4 bullets -> heading + 6 bullets
# ---------------------------------------------------------------------
Replace-Text "Political Field Operations and Data Management" "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns"

Replace-Text "• Managed field operations for political campaigns and research projects" "• Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions"

Replace-Text "• Developed data collection and management systems for field work" "• Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm"

Replace-Text "• Trained field staff on data collection protocols and quality control" "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"

Replace-Text "• Analyzed field data to inform campaign strategy and research findings" "• Created custom reports and data visualizations based on specific client requirements"

$p3 = Add-BulletAfter "Created custom reports and data visualizations based on specific client requirements" "• Introduced mapping and geospatial analysis into standard reporting procedures"
if ($p3) {
    $p3.Range.InsertParagraphAfter()
    $p4 = $p3.Next()
    $p4.Range.Text = "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL"
}

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
